$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date column (E) entries from 2-digit year (DD.MM.YY) to 4-digit year (DD.MM.YYYY)
$ws.Range("E2").Value = "14.03.2023"
$ws.Range("E3").Value = "16.03.2023"
$ws.Range("E4").Value = "21.03.2023"
$ws.Range("E5").Value = "22.03.2023"
$ws.Range("E6").Value = "24.03.2023"

# Update the active selection to G5, as reflected in the sheet view
$ws.Range("G5").Select()
